# Corrected the population and the GDP of Slovenia (row 19) to match csv/excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slovenia population figure (column C) was wrong by a factor of 10 -> fix it.
# The GDP column (B19) is a formula (" 34400 * C19") so it recalculates
# automatically once the population is corrected.
$ws.Range("C19").Value = 2239490

# Leave the selection on the cell that was actually edited/reviewed.
$ws.Range("B19").Select()
